# Applies the "Updated cryptos list" data refresh to Sheet1.
# Source cells are t="inlineStr" text cells (prices/links/percent deltas),
# so plain numeric-looking strings are written through Set-TextValue to
# stop Excel from auto-converting them to the Number type (and, e.g.,
# dropping the trailing zero in "55.10").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "37.137.32"
# Row 3
$ws.Range("D3").Value = "2.023.97"
$ws.Range("E3").Value = "  -2.99%  "
# Row 4
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
Set-TextValue $ws.Range("D5") "227.27"
$ws.Range("E5").Value = "  -2.75%  "
# Row 6
Set-TextValue $ws.Range("D6") "0.609"
$ws.Range("E6").Value = "  -4.56%  "
# Row 7
$ws.Range("E7").Value = "  +0.05%  "
# Row 8
Set-TextValue $ws.Range("D8") "55.10"
$ws.Range("E8").Value = "  -5.12%  "
# Row 9
$ws.Range("E9").Value = "  -2.95%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.0789"
$ws.Range("E10").Value = "  +1.22%  "
# Row 11
$ws.Range("E11").Value = "  -3.55%  "
# Row 12
$ws.Range("D12").Value = "2.324.70"
$ws.Range("E12").Value = "  -2.93%  "
# Row 13
Set-TextValue $ws.Range("D13") "14.28"
$ws.Range("E13").Value = "  -5.92%  "
# Row 14
Set-TextValue $ws.Range("D14") "20.42"
$ws.Range("E14").Value = "  -3.17%  "
# Row 15
Set-TextValue $ws.Range("D15") "0.744"
$ws.Range("E15").Value = "  -4.23%  "
# Row 16
$ws.Range("E16").Value = "  -3.23%  "
# Row 17
$ws.Range("D17").Value = "2.023.19"
$ws.Range("E17").Value = "  -2.80%  "
# Row 18
$ws.Range("D18").Value = "37.008.86"
$ws.Range("E18").Value = "  -1.88%  "
# Row 19
Set-TextValue $ws.Range("D19") "6.03"
$ws.Range("E19").Value = "  -1.30%  "
# Row 20
Set-TextValue $ws.Range("D20") "68.83"
# Row 21
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  -0.13%  "
# Row 22
Set-TextValue $ws.Range("D22") "223.30"
$ws.Range("E22").Value = "  -2.65%  "
# Row 23
$ws.Range("E23").Value = "  +0.25%  "
# Row 24
Set-TextValue $ws.Range("D24") "2.40"
$ws.Range("E24").Value = "  +0.77%  "
# Row 25
Set-TextValue $ws.Range("D25") "2.26"
$ws.Range("E25").Value = "  -5.28%  "
# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "9.37"
$ws.Range("E26").Value = "  -3.70%  "
# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D27") "167.84"
$ws.Range("E27").Value = "  -1.84%  "
# Row 28
$ws.Range("E28").Value = "  -6.71%  "
# Row 29
Set-TextValue $ws.Range("D29") "18.75"
$ws.Range("E29").Value = "  -3.89%  "
# Row 30
$ws.Range("E30").Value = "  -4.22%  "
# Row 31
$ws.Range("E31").Value = "  -5.03%  "
# Row 32
$ws.Range("E32").Value = "  -4.31%  "
# Row 33
Set-TextValue $ws.Range("D33") "0.0606"
$ws.Range("E33").Value = "  -4.54%  "
# Row 34
Set-TextValue $ws.Range("D34") "4.46"
$ws.Range("E34").Value = "  -2.91%  "
# Row 35
Set-TextValue $ws.Range("D35") "2.36"
$ws.Range("E35").Value = "  -5.37%  "
# Row 36
Set-TextValue $ws.Range("D36") "1.82"
$ws.Range("E36").Value = "  -0.10%  "
# Row 37
$ws.Range("E37").Value = "  +0.12%  "
# Row 38
$ws.Range("E38").Value = "  -5.14%  "
# Row 39
Set-TextValue $ws.Range("D39") "5.36"
$ws.Range("E39").Value = "  -0.22%  "
# Row 40
$ws.Range("D40").Value = "1.493.75"
$ws.Range("E40").Value = "  +2.82%  "
# Row 41
$ws.Range("E41").Value = "  -7.34%  "
# Row 42
$ws.Range("E42").Value = "  -2.08%  "
# Row 43
$ws.Range("E43").Value = "  -4.02%  "
# Row 44
Set-TextValue $ws.Range("D44") "94.94"
$ws.Range("E44").Value = "  -6.04%  "
# Row 45
Set-TextValue $ws.Range("D45") "16.51"
$ws.Range("E45").Value = "  -0.74%  "
# Row 46
$ws.Range("E46").Value = "  -5.40%  "
# Row 47
$ws.Range("E47").Value = "  -5.00%  "
# Row 48
$ws.Range("E48").Value = "  -1.07%  "
# Row 49
$ws.Range("E49").Value = "  -1.61%  "
# Row 50
Set-TextValue $ws.Range("D50") "3.67"
$ws.Range("E50").Value = "  -10.84%  "
# Row 51
$ws.Range("D51").Value = "2.214.13"
$ws.Range("E51").Value = "  -2.86%  "
